# Insert one new data row at row 789 on Sheet1 ("sei2"), shifting the
# existing rows 789:830 down to 790:831, and populate the new row with
# the 2026/02/08 (Sunday) / 17:00 reading, ranking 201 — the entry that
# was missing before this "daily auto push" commit.
#
# Everything from the former row 789 onward keeps its values; only the
# physical row position moves down by one (dimension grows from
# A1:D830 to A1:D831).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 789..830 down by one row, leaving row 789 blank.
$ws.Rows.Item(789).Insert()

# Column A holds dates stored as plain text (e.g. "2026/02/08"), not
# Excel date serials, just like every other row in the sheet. Force
# text formatting before the assignment so COM doesn't auto-coerce the
# date-shaped string into a date value, then drop back to the
# worksheet's default (unstyled) cell style to match the rest of the
# data rows, which carry no explicit style.
$newDateCell = $ws.Range("A789")
$newDateCell.NumberFormat = "@"
$newDateCell.Value = "2026/02/08"
$newDateCell.Style = "Normal"

$ws.Range("B789").Value = "日"
$ws.Range("C789").Value = 17
$ws.Range("D789").Value = 201
